$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.277.50'
$ws.Range("E2").Value = '  -2.28%  '

$ws.Range("D3").Value = '1.863.59'
$ws.Range("E3").Value = '  -2.01%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.56'
$ws.Range("E5").Value = '  -1.45%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  -0.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4390'
$ws.Range("E7").Value = '  -4.51%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3739'
$ws.Range("E8").Value = '  -1.89%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07546'
$ws.Range("E9").Value = '  -2.03%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9391'
$ws.Range("E10").Value = '  -3.54%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.30'
$ws.Range("E11").Value = '  -2.91%  '

$ws.Range("D12").Value = '1.868.43'
$ws.Range("E12").Value = '  -2.20%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.728'
$ws.Range("E13").Value = '  -2.76%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.478'
$ws.Range("E14").Value = '  -2.99%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06860'
$ws.Range("E15").Value = '  -2.95%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '82.10'
$ws.Range("E17").Value = '  -2.12%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009104'

$ws.Range("E19").Value = '  -0.23%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.03'
$ws.Range("E20").Value = '  -3.44%  '

$ws.Range("D21").Value = '28.279.41'
$ws.Range("E21").Value = '  -2.16%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.153'
$ws.Range("E22").Value = '  -2.31%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.75'
$ws.Range("E23").Value = '  -0.91%  '

$ws.Range("D24").Value = '2.113.76'
$ws.Range("E24").Value = '  -1.24%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.041'
$ws.Range("E25").Value = '  -2.65%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.86'
$ws.Range("E26").Value = '  -1.99%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.45'
$ws.Range("E27").Value = '  -2.86%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.367'
$ws.Range("E28").Value = '  -4.48%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.63'
$ws.Range("E29").Value = '  -2.37%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.734'
$ws.Range("E30").Value = '  -5.56%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09058'
$ws.Range("E31").Value = '  -2.09%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.8063'
$ws.Range("E32").Value = '  -5.88%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.861'
$ws.Range("E33").Value = '  -4.20%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.175'
$ws.Range("E34").Value = '  -5.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.943'
$ws.Range("E35").Value = '  +0.12%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.003'
$ws.Range("E36").Value = '  -0.12%  '

$ws.Range("E37").Value = '  -0.95%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05476'
$ws.Range("E38").Value = '  -3.48%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.008'
$ws.Range("E39").Value = '  +9.17%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01960'
$ws.Range("E40").Value = '  -3.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.178'
$ws.Range("E41").Value = '  -2.71%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5260'
$ws.Range("E42").Value = '  -3.96%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1678'
$ws.Range("E43").Value = '  -4.23%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.829'
$ws.Range("E44").Value = '  -4.78%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.075'
$ws.Range("E45").Value = '  +0.89%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06782'
$ws.Range("E46").Value = '  -0.55%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4887'
$ws.Range("E47").Value = '  -5.17%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000002525'
$ws.Range("E48").Value = '  -2.08%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.56'
$ws.Range("E49").Value = '  -5.87%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '107.62'
$ws.Range("E50").Value = '  -2.27%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.686'
$ws.Range("E51").Value = '  -4.48%  '
